$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.06842864744397358
$ws.Range("J2").Value = 0.06842864744397358
$ws.Range("M2").Value = 10.50827633333333
$ws.Range("N2").Value = 31.524829
$ws.Range("O2").Value = 0.1682660991018133
$ws.Range("P2").Value = 0.1682660991018134
$ws.Range("Q2").Value = 0.8036064215565556
$ws.Range("R2").Value = 7.232457794009001
$ws.Range("S2").Value = 0.0115142215722107
$ws.Range("T2").Value = 0.01151422157221071
$ws.Range("I3").Value = 0.06842864744397358
$ws.Range("J3").Value = 0.06842864744397358
$ws.Range("O3").Value = 0.4955285863849104
$ws.Range("P3").Value = 0.4955285863849105
$ws.Range("S3").Value = 0.03390835093614364
$ws.Range("T3").Value = 0.03390835093614365
$ws.Range("I4").Value = 0.06842864744397358
$ws.Range("J4").Value = 0.06842864744397358
$ws.Range("M4").Value = 6.495209666666667
$ws.Range("N4").Value = 19.485629
$ws.Range("O4").Value = 0.1040059814559238
$ws.Range("P4").Value = 0.1040059814559238
$ws.Range("Q4").Value = 0.4967124989787778
$ws.Range("R4").Value = 4.470412490809
$ws.Range("S4").Value = 0.007116988637111862
$ws.Range("T4").Value = 0.007116988637111862
$ws.Range("I5").Value = 0.06842864744397358
$ws.Range("J5").Value = 0.06842864744397358
$ws.Range("M5").Value = 9.909791666666667
$ws.Range("N5").Value = 29.729375
$ws.Range("O5").Value = 0.1586827309986352
$ws.Range("P5").Value = 0.1586827309986352
$ws.Range("Q5").Value = 0.7578381046527779
$ws.Range("R5").Value = 6.820542941875001
$ws.Range("S5").Value = 0.0108584446549525
$ws.Range("T5").Value = 0.0108584446549525
$ws.Range("I6").Value = 0.06842864744397358
$ws.Range("J6").Value = 0.06842864744397358
$ws.Range("M6").Value = 4.591137333333333
$ws.Range("N6").Value = 13.773412
$ws.Range("O6").Value = 0.07351660205871713
$ws.Range("P6").Value = 0.07351660205871713
$ws.Range("Q6").Value = 0.3511011060502223
$ws.Range("R6").Value = 3.159909954452
$ws.Range("S6").Value = 0.005030641643554857
$ws.Range("T6").Value = 0.005030641643554857
$ws.Range("G7").Value = 0.621785
$ws.Range("H7").Value = 1.865355
$ws.Range("I7").Value = 0.5563733034589394
$ws.Range("J7").Value = 0.5563733034589394
$ws.Range("M7").Value = 10.50827633333333
$ws.Range("N7").Value = 31.524829
$ws.Range("O7").Value = 0.1682660991018133
$ws.Range("P7").Value = 0.1682660991018134
$ws.Range("Q7").Value = 6.533888599921667
$ws.Range("R7").Value = 58.804997399295
$ws.Range("S7").Value = 0.09361876541742517
$ws.Range("T7").Value = 0.09361876541742518
$ws.Range("G8").Value = 0.621785
$ws.Range("H8").Value = 1.865355
$ws.Range("I8").Value = 0.5563733034589394
$ws.Range("J8").Value = 0.5563733034589394
$ws.Range("O8").Value = 0.4955285863849104
$ws.Range("P8").Value = 0.4955285863849105
$ws.Range("Q8").Value = 19.24171653588167
$ws.Range("R8").Value = 173.175448822935
$ws.Range("S8").Value = 0.2756988765653111
$ws.Range("T8").Value = 0.2756988765653111
$ws.Range("G9").Value = 0.621785
$ws.Range("H9").Value = 1.865355
$ws.Range("I9").Value = 0.5563733034589394
$ws.Range("J9").Value = 0.5563733034589394
$ws.Range("M9").Value = 6.495209666666667
$ws.Range("N9").Value = 19.485629
$ws.Range("O9").Value = 0.1040059814559238
$ws.Range("P9").Value = 0.1040059814559238
$ws.Range("Q9").Value = 4.038623942588334
$ws.Range("R9").Value = 36.347615483295
$ws.Range("S9").Value = 0.05786615148212151
$ws.Range("T9").Value = 0.05786615148212151
$ws.Range("G10").Value = 0.621785
$ws.Range("H10").Value = 1.865355
$ws.Range("I10").Value = 0.5563733034589394
$ws.Range("J10").Value = 0.5563733034589394
$ws.Range("M10").Value = 9.909791666666667
$ws.Range("N10").Value = 29.729375
$ws.Range("O10").Value = 0.1586827309986352
$ws.Range("P10").Value = 0.1586827309986352
$ws.Range("Q10").Value = 6.161759811458334
$ws.Range("R10").Value = 55.45583830312501
$ws.Range("S10").Value = 0.0882868352475969
$ws.Range("T10").Value = 0.0882868352475969
$ws.Range("G11").Value = 0.621785
$ws.Range("H11").Value = 1.865355
$ws.Range("I11").Value = 0.5563733034589394
$ws.Range("J11").Value = 0.5563733034589394
$ws.Range("M11").Value = 4.591137333333333
$ws.Range("N11").Value = 13.773412
$ws.Range("O11").Value = 0.07351660205871713
$ws.Range("P11").Value = 0.07351660205871713
$ws.Range("Q11").Value = 2.854700326806667
$ws.Range("R11").Value = 25.69230294126
$ws.Range("S11").Value = 0.04090267474648471
$ws.Range("T11").Value = 0.04090267474648471
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.4193093333333334
$ws.Range("H12").Value = 1.257928
$ws.Range("I12").Value = 0.375198049097087
$ws.Range("J12").Value = 0.375198049097087
$ws.Range("M12").Value = 10.50827633333333
$ws.Range("N12").Value = 31.524829
$ws.Range("O12").Value = 0.1682660991018133
$ws.Range("P12").Value = 0.1682660991018134
$ws.Range("Q12").Value = 4.406218343812444
$ws.Range("R12").Value = 39.65596509431201
$ws.Range("S12").Value = 0.06313311211217747
$ws.Range("T12").Value = 0.06313311211217748
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.4193093333333334
$ws.Range("H13").Value = 1.257928
$ws.Range("I13").Value = 0.375198049097087
$ws.Range("J13").Value = 0.375198049097087
$ws.Range("O13").Value = 0.4955285863849104
$ws.Range("P13").Value = 0.4955285863849105
$ws.Range("Q13").Value = 12.97591825606844
$ws.Range("R13").Value = 116.783264304616
$ws.Range("S13").Value = 0.1859213588834557
$ws.Range("T13").Value = 0.1859213588834558
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.4193093333333334
$ws.Range("H14").Value = 1.257928
$ws.Range("I14").Value = 0.375198049097087
$ws.Range("J14").Value = 0.375198049097087
$ws.Range("M14").Value = 6.495209666666667
$ws.Range("N14").Value = 19.485629
$ws.Range("O14").Value = 0.1040059814559238
$ws.Range("P14").Value = 0.1040059814559238
$ws.Range("Q14").Value = 2.723502035190223
$ws.Range("R14").Value = 24.511518316712
$ws.Range("S14").Value = 0.03902284133669041
$ws.Range("T14").Value = 0.03902284133669041
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.4193093333333334
$ws.Range("H15").Value = 1.257928
$ws.Range("I15").Value = 0.375198049097087
$ws.Range("J15").Value = 0.375198049097087
$ws.Range("M15").Value = 9.909791666666667
$ws.Range("N15").Value = 29.729375
$ws.Range("O15").Value = 0.1586827309986352
$ws.Range("P15").Value = 0.1586827309986352
$ws.Range("Q15").Value = 4.155268137222222
$ws.Range("R15").Value = 37.39741323500001
$ws.Range("S15").Value = 0.05953745109608576
$ws.Range("T15").Value = 0.05953745109608578
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.4193093333333334
$ws.Range("H16").Value = 1.257928
$ws.Range("I16").Value = 0.375198049097087
$ws.Range("J16").Value = 0.375198049097087
$ws.Range("M16").Value = 4.591137333333333
$ws.Range("N16").Value = 13.773412
$ws.Range("O16").Value = 0.07351660205871713
$ws.Range("P16").Value = 0.07351660205871713
$ws.Range("Q16").Value = 1.925106734481778
$ws.Range("R16").Value = 17.325960610336
$ws.Range("S16").Value = 0.02758328566867756
$ws.Range("T16").Value = 0.02758328566867756
